# Append a new Lancers work-listing row (2026-01-02 01:59:58 JST) and
# refresh the "retrieved at" timestamp on the most recent existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Refresh row 2's timestamp to the latest fetch time.
$ws.Range("A2").Value = "2026-01-02 01:59:58"

# Append the newly scraped listing as row 3.
$ws.Range("A3").Value = "2026-01-02 01:59:58"
$ws.Range("B3").Value = "進行管理およびチームディレクションを担当"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "~ 5,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = "◇管理"

# Hyperlink the URL cell, matching the style used for the existing link.
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5418064")
$ws.Range("F3").Style = "Hyperlink"
